$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell values, entered in the same order the author typed them (this is
# what determines the shared-string table order) ---
$ws.Range("B1").Value = "UserName"
$ws.Range("C1").Value = "Password"
$ws.Range("D1").Value = "FirstName"
$ws.Range("E1").Value = "LastName"
$ws.Range("F1").Value = "Email"

$ws.Range("A2").Value = "AccountDetails"
$ws.Range("B2").Value = "sbingi@helenoftroy.com"
$ws.Range("C2").Value = "Auislzkuakm03!"

$ws.Range("G1").Value = "CustomerID"

$ws.Range("A3").Value = "ProDeal"
$ws.Range("G3").Value = "'76"

$ws.Range("H1").Value = "AssociationEmail"
$ws.Range("H3").Value = "skasarla@helenoftroy.com"

$ws.Range("D4").Value = "QA"
$ws.Range("E4").Value = "TEST"

$ws.Range("A4").Value = "ProDealForm"
$ws.Range("H4").Value = "nsada@helenoftroy.com"

$ws.Range("I1").Value = "Association"
$ws.Range("I4").Value = "Testing"

$ws.Range("J1").Value = "Status"
$ws.Range("J4").Value = "Pending"

$ws.Range("K1").Value = "Group"
$ws.Range("L1").Value = "Comment"

$ws.Range("K4").Value = "General"
$ws.Range("L4").Value = "Qa testing the prodeal form"

$ws.Range("G4").Value = 976
